$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (ExpectedFilenames) values need updating: a new "StandardExcelReport...-2023_"
# row is inserted ahead of each existing Excel/Word report pair, and the existing
# "ExcelReport-" naming loses the space around the dash after NewImportLogic_1.
$ws.Cells.Item(2, 7).Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-2023_"
$ws.Cells.Item(3, 7).Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Clinical-"
$ws.Cells.Item(4, 7).Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"

$ws.Cells.Item(5, 7).Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Economic-2023_"
$ws.Cells.Item(6, 7).Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Economic-"
$ws.Cells.Item(7, 7).Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Economic-"

$ws.Cells.Item(8, 7).Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-2023_"
$ws.Cells.Item(9, 7).Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Quality of Life-"
$ws.Cells.Item(10, 7).Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-"

$ws.Cells.Item(11, 7).Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-2023_"
$ws.Cells.Item(12, 7).Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Real-world Evidence-"
$ws.Cells.Item(13, 7).Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-"

# Make column G fit its new (longer) content, matching the widened column in the workbook.
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null

# Restore the selection to the cell the author last had selected.
$ws.Range("G13").Select() | Out-Null
